# The "Förändrad" (Changed) column C holds a date serial that is bumped by
# one day for every data row (rows 2-255) on each automatic refresh of the
# sheet. Increment each C-column cell from 46061 to 46062.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C255")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
